$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay formatted as Text so values like
# "69.560.63" or "1.00" are not coerced to numbers and lose formatting.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '69.560.63'
$ws.Range("E2").Value = '  +0.22%  '

# Row 3
$ws.Range("D3").Value = '3.692.67'
$ws.Range("E3").Value = '  +0.12%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").Value = '676.81'
$ws.Range("E5").Value = '  -1.34%  '

# Row 6
$ws.Range("D6").Value = '161.59'
$ws.Range("E6").Value = '  +0.97%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  +0.27%  '

# Row 9
$ws.Range("E9").Value = '  +0.75%  '

# Row 10
$ws.Range("E10").Value = '  -0.44%  '

# Row 11
$ws.Range("D11").Value = '0.441'
$ws.Range("E11").Value = '  +1.26%  '

# Row 12
$ws.Range("D12").Value = '0.0000234'
$ws.Range("E12").Value = '  +0.32%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '32.50'
$ws.Range("E13").Value = '  +0.14%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '3.674.60'
$ws.Range("E14").Value = '  -0.21%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '69.470.84'
$ws.Range("E15").Value = '  +0.04%  '

# Row 16
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").Value = '0.117'
$ws.Range("E16").Value = '  +2.35%  '

# Row 17
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '16.01'
$ws.Range("E17").Value = '  +0.49%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '6.48'
$ws.Range("E18").Value = '  +0.08%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '471.11'
$ws.Range("E19").Value = '  +0.79%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '9.80'
$ws.Range("E20").Value = '  -2.09%  '

# Row 21
$ws.Range("B21").Value = 'Polygon'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D21").Value = '0.650'
$ws.Range("E21").Value = '  +0.41%  '

# Row 22
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").Value = '80.54'
$ws.Range("E22").Value = '  +1.38%  '

# Row 23
$ws.Range("B23").Value = 'WrappedeETH'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D23").Value = '3.837.71'
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.14%  '

# Row 25
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '0.0000126'
$ws.Range("E25").Value = '  +1.18%  '

# Row 26
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '10.86'
$ws.Range("E26").Value = '  -1.15%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '9.13'
$ws.Range("E27").Value = '  -0.61%  '

# Row 28
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.70'
$ws.Range("E28").Value = '  -0.53%  '

# Row 29
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").Value = '1.74'
$ws.Range("E29").Value = '  -0.10%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '2.02'
$ws.Range("E30").Value = '  -0.07%  '

# Row 31
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '6.59'
$ws.Range("E31").Value = '  -0.53%  '

# Row 32
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.39%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '27.00'
$ws.Range("E33").Value = '  +0.93%  '

# Row 34
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '3.680.82'
$ws.Range("E34").Value = '  +0.46%  '

# Row 35
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.162'
$ws.Range("E35").Value = '  +0.72%  '

# Row 36
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '8.47'
$ws.Range("E36").Value = '  +3.61%  '

# Row 37
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = '6.23'
$ws.Range("E37").Value = '  +0.83%  '

# Row 38
$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.00%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '2.26'
$ws.Range("E39").Value = '  -1.20%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.10%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.0901'
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '168.15'
$ws.Range("E42").Value = '  +0.56%  '

# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '0.943'
$ws.Range("E43").Value = '  -0.01%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '46.54'
$ws.Range("E44").Value = '  -2.91%  '

# Row 45
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.75'
$ws.Range("E45").Value = '  +0.57%  '

# Row 46
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Value = '0.000278'
$ws.Range("E46").Value = '  +0.82%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '27.98'
$ws.Range("E47").Value = '  -2.68%  '

# Row 48
$ws.Range("D48").Value = '1.29'
$ws.Range("E48").Value = '  -1.30%  '

# Row 49
$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").Value = '1.09'
$ws.Range("E49").Value = '  -2.47%  '

# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '7.89'
$ws.Range("E50").Value = '  +1.16%  '

# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '0.266'
$ws.Range("E51").Value = '  +1.48%  '
